$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2736
$ws.Range("J17").Value = 2736
$ws.Range("L17").Value = 8208
$ws.Range("N17").Value = -8544

$ws.Range("H132").Value = 5014.4136
$ws.Range("I132").Value = 4780.72
$ws.Range("K132").Value = 14342.16
$ws.Range("M132").Value = -11812.16

$ws.Range("H138").Value = 2284.6902
$ws.Range("I138").Value = 842.129
$ws.Range("J138").Value = 3402.675
$ws.Range("K138").Value = 2526.387
$ws.Range("L138").Value = 10208.025
$ws.Range("M138").Value = 2613.613
$ws.Range("N138").Value = -20488.025

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H32").Value = 5229.1577
$ws.Range("I32").Value = 5055.8193
$ws.Range("K32").Value = 5055.8193
$ws.Range("M32").Value = -4768.8193

$ws.Range("H45").Value = 100402.336
$ws.Range("I45").Value = 187107.19
$ws.Range("K45").Value = 187107.19
$ws.Range("M45").Value = -186730.19

$ws.Range("H106").Value = 60000
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H132").Value = 2528.5938
$ws.Range("I132").Value = 2193.625
$ws.Range("J132").Value = 3533.5
$ws.Range("K132").Value = 6580.875
$ws.Range("L132").Value = 10600.5
$ws.Range("M132").Value = -4050.875
$ws.Range("N132").Value = -15660.5

$ws.Range("H135").Value = 87475.664
$ws.Range("J135").Value = 87475.664
$ws.Range("L135").Value = 87475.664
$ws.Range("N135").Value = -97615.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H99").Value = 13361.04
$ws.Range("I99").Value = 15678.944
$ws.Range("K99").Value = 15678.944
$ws.Range("M99").Value = -14180.944

$ws.Range("H134").Value = 6161.04
$ws.Range("I134").Value = 6325.1665
$ws.Range("K134").Value = 18975.4995
$ws.Range("M134").Value = -16440.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12158.4
$ws.Range("I31").Value = 16264.889
$ws.Range("J31").Value = 5998.6665
$ws.Range("K31").Value = 16264.889
$ws.Range("L31").Value = 5998.6665
$ws.Range("M31").Value = -15969.889
$ws.Range("N31").Value = -6588.6665

$ws.Range("H34").Value = 12158.4
$ws.Range("I34").Value = 16264.889
$ws.Range("J34").Value = 5998.6665
$ws.Range("K34").Value = 16264.889
$ws.Range("L34").Value = 5998.6665
$ws.Range("M34").Value = -16062.889
$ws.Range("N34").Value = -6402.6665

$ws.Range("H58").Value = 4794.2856
$ws.Range("I58").Value = 4970.4443
$ws.Range("K58").Value = 4970.4443
$ws.Range("M58").Value = -4767.4443

$ws.Range("H86").Value = 7129.294
$ws.Range("I86").Value = 6798.875
$ws.Range("J86").Value = 7423
$ws.Range("K86").Value = 6798.875
$ws.Range("L86").Value = 7423
$ws.Range("M86").Value = -5675.875
$ws.Range("N86").Value = -9669

$ws.Range("H89").Value = 7129.294
$ws.Range("I89").Value = 6798.875
$ws.Range("J89").Value = 7423
$ws.Range("K89").Value = 33994.375
$ws.Range("L89").Value = 37115
$ws.Range("M89").Value = -28378.375
$ws.Range("N89").Value = -48347

$ws.Range("H107").Value = 4267.706
$ws.Range("I107").Value = 6615.85
$ws.Range("K107").Value = 6615.85
$ws.Range("M107").Value = -4695.85

$ws.Range("H134").Value = 2509.6206
$ws.Range("I134").Value = 2631.4546
$ws.Range("J134").Value = 2126.7144
$ws.Range("K134").Value = 7894.3638
$ws.Range("L134").Value = 6380.1432
$ws.Range("M134").Value = -5359.3638
$ws.Range("N134").Value = -11450.1432

$ws.Range("H136").Value = 4794.2856
$ws.Range("I136").Value = 4970.4443
$ws.Range("K136").Value = 14911.3329
$ws.Range("M136").Value = -12361.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 44521636
$ws.Range("J4").Value = 105429240
$ws.Range("L4").Value = 316287720
$ws.Range("N4").Value = -316287944

$ws.Range("H12").Value = 107
$ws.Range("I12").Value = 248
$ws.Range("K12").Value = 744
$ws.Range("M12").Value = -571

$ws.Range("H107").Value = 992.4643
$ws.Range("J107").Value = 1231.95
$ws.Range("L107").Value = 3695.85
$ws.Range("N107").Value = -7535.85

$ws.Range("H113").Value = 1118.5714
$ws.Range("J113").Value = 1225.8334
$ws.Range("L113").Value = 3677.5002
$ws.Range("N113").Value = -8017.5002

$ws.Range("H131").Value = 1594.8395
$ws.Range("I131").Value = 1299.7142
$ws.Range("J131").Value = 1622.7567
$ws.Range("K131").Value = 3899.1426
$ws.Range("L131").Value = 4868.2701
$ws.Range("M131").Value = 1140.8574
$ws.Range("N131").Value = -14948.2701

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 67156570
$ws.Range("I2").Value = 170006400
$ws.Range("K2").Value = 170006400
$ws.Range("M2").Value = -170006288

$ws.Range("H46").Value = 5362620.5
$ws.Range("J46").Value = 8713664
$ws.Range("L46").Value = 8713664
$ws.Range("N46").Value = -8714040

$ws.Range("H135").Value = 94226.87
$ws.Range("J135").Value = 94226.87
$ws.Range("L135").Value = 94226.87
$ws.Range("N135").Value = -104366.87

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 149834
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 149834
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H81").Value = 15261.875
$ws.Range("I81").Value = 26026.25
$ws.Range("K81").Value = 52052.5
$ws.Range("M81").Value = -50991.5

$ws.Range("H84").Value = 15261.875
$ws.Range("I84").Value = 26026.25
$ws.Range("K84").Value = 260262.5
$ws.Range("M84").Value = -254958.5

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H122").Value = 5258.275
$ws.Range("I122").Value = 2664.3125
$ws.Range("K122").Value = 7992.9375
$ws.Range("M122").Value = -5542.9375

$ws.Range("H126").Value = 29557.732
$ws.Range("I126").Value = 32567
$ws.Range("K126").Value = 97701
$ws.Range("M126").Value = -95231
